# Test_Cases.xlsx - add new test-case rows for spin 3 (rows 23-55)
# and requirement-id updates, per the commit: "Test results and updated
# spreadsheets for new tests in spin 3. Updated test descriptions with
# requirement ids."
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 23: length_delimited_12_02
$ws.Range("B23").Value = 'length_delimited_12_02'
$ws.Range("C23").Value = 'Passed'
$ws.Range("D23").Value = 'DFDL-12-043R'
$ws.Range("E23").Value = 'High'
$ws.Range("F23").Value = 2
$ws.Range("G23").Value = 'text'
$ws.Range("H23").Value = 'DFDL-112'

# Row 24: length_delimited_12_03
$ws.Range("B24").Value = 'length_delimited_12_03'
$ws.Range("C24").Value = 'Passed'
$ws.Range("D24").Value = 'DFDL-12-043R'
$ws.Range("E24").Value = 'High'
$ws.Range("F24").Value = 2
$ws.Range("G24").Value = 'text'
$ws.Range("H24").Value = 'DFDL-112'

# Row 25: introduction_1_02
$ws.Range("B25").Value = 'introduction_1_02'
$ws.Range("C25").Value = 'Passed'
$ws.Range("D25").Value = 'DFDL-12-043R'
$ws.Range("E25").Value = 'High'
$ws.Range("F25").Value = 2
$ws.Range("G25").Value = 'text'
$ws.Range("H25").Value = 'DFDL-112'

# Row 26: multiple_delimiters
$ws.Range("B26").Value = 'multiple_delimiters'
$ws.Range("C26").Value = 'Passed'
$ws.Range("D26").Value = 'DFDL-12-033R'
$ws.Range("E26").Value = 'High'
$ws.Range("F26").Value = 2
$ws.Range("G26").Value = 'text'
$ws.Range("H26").Value = 'DFDL-112'

# Row 27: NumSeq_01
$ws.Range("B27").Value = 'NumSeq_01'
$ws.Range("C27").Value = 'Passed'
$ws.Range("D27").Value = 'DFDL-12-048R'
$ws.Range("E27").Value = 'High'
$ws.Range("F27").Value = 2
$ws.Range("G27").Value = 'text'
$ws.Range("H27").Value = 'DFDL-112'

# Row 28: NumSeq_03
$ws.Range("B28").Value = 'NumSeq_03'
$ws.Range("C28").Value = 'Passed'
$ws.Range("D28").Value = 'DFDL-12-045R'
$ws.Range("E28").Value = 'High'
$ws.Range("F28").Value = 2
$ws.Range("G28").Value = 'text'
$ws.Range("H28").Value = 'DFDL-112'

# Row 29: NumSeq_04
$ws.Range("B29").Value = 'NumSeq_04'
$ws.Range("C29").Value = 'Passed'
$ws.Range("D29").Value = 'DFDL-12-042R'
$ws.Range("E29").Value = 'High'
$ws.Range("F29").Value = 2
$ws.Range("G29").Value = 'text'
$ws.Range("H29").Value = 'DFDL-112'

# Row 30: AB000
$ws.Range("B30").Value = 'AB000'
$ws.Range("C30").Value = 'Passed'
$ws.Range("D30").Value = 'DFDL-12-043R'
$ws.Range("E30").Value = 'High'
$ws.Range("F30").Value = 2
$ws.Range("G30").Value = 'text'
$ws.Range("H30").Value = 'DFDL-112'

# Row 31: AB001
$ws.Range("B31").Value = 'AB001'
$ws.Range("C31").Value = 'Passed'
$ws.Range("D31").Value = 'DFDL-12-043R'
$ws.Range("E31").Value = 'High'
$ws.Range("F31").Value = 2
$ws.Range("G31").Value = 'text'
$ws.Range("H31").Value = 'DFDL-112'

# Row 32: AB002
$ws.Range("B32").Value = 'AB002'
$ws.Range("C32").Value = 'Passed'
$ws.Range("D32").Value = 'DFDL-12-043R'
$ws.Range("E32").Value = 'High'
$ws.Range("F32").Value = 2
$ws.Range("G32").Value = 'text'
$ws.Range("H32").Value = 'DFDL-112'

# Row 33: AB003
$ws.Range("B33").Value = 'AB003'
$ws.Range("C33").Value = 'Passed'
$ws.Range("D33").Value = 'DFDL-12-043R'
$ws.Range("E33").Value = 'High'
$ws.Range("F33").Value = 2
$ws.Range("G33").Value = 'text'
$ws.Range("H33").Value = 'DFDL-112'

# Row 34: AN000
$ws.Range("B34").Value = 'AN000'
$ws.Range("C34").Value = 'Passed'
$ws.Range("D34").Value = 'DFDL-12-042R'
$ws.Range("E34").Value = 'High'
$ws.Range("F34").Value = 2
$ws.Range("G34").Value = 'text'
$ws.Range("H34").Value = 'DFDL-112'

# Row 35: AN001
$ws.Range("B35").Value = 'AN001'
$ws.Range("C35").Value = 'Passed'
$ws.Range("D35").Value = 'DFDL-12-042R'
$ws.Range("E35").Value = 'High'
$ws.Range("F35").Value = 2
$ws.Range("G35").Value = 'text'
$ws.Range("H35").Value = 'DFDL-112'

# Row 36: AI000_rev
$ws.Range("B36").Value = 'AI000_rev'
$ws.Range("C36").Value = 'Passed'
$ws.Range("D36").Value = 'DFDL-12-087R'
$ws.Range("E36").Value = 'High'
$ws.Range("F36").Value = 3
$ws.Range("G36").Value = 'text'
$ws.Range("H36").Value = 'DFDL-205'

# Row 37: LengthKindPattern
$ws.Range("B37").Value = 'LengthKindPattern'
$ws.Range("C37").Value = 'Passed'
$ws.Range("D37").Value = 'DFDL-12-088R'
$ws.Range("E37").Value = 'High'
$ws.Range("F37").Value = 3
$ws.Range("G37").Value = 'text'
$ws.Range("H37").Value = 'DFDL-207'

# Row 38: LengthKindPatternCompound
$ws.Range("B38").Value = 'LengthKindPatternCompound'
$ws.Range("C38").Value = 'Passed'
$ws.Range("D38").Value = 'DFDL-12-088R'
$ws.Range("E38").Value = 'High'
$ws.Range("F38").Value = 3
$ws.Range("G38").Value = 'text'
$ws.Range("H38").Value = 'DFDL-205'

# Row 39: lengthKindPattern_01
$ws.Range("B39").Value = 'lengthKindPattern_01'
$ws.Range("C39").Value = 'Passed'
$ws.Range("D39").Value = 'DFDL-12-088R'
$ws.Range("E39").Value = 'High'
$ws.Range("F39").Value = 3
$ws.Range("G39").Value = 'text'
$ws.Range("H39").Value = 'DFDL-205'

# Row 40: lengthKindPattern_02
$ws.Range("B40").Value = 'lengthKindPattern_02'
$ws.Range("C40").Value = 'Passed'
$ws.Range("D40").Value = 'DFDL-12-088R'
$ws.Range("E40").Value = 'High'
$ws.Range("F40").Value = 3
$ws.Range("G40").Value = 'text'
$ws.Range("H40").Value = 'DFDL-205'

# Row 41: lengthKindPattern_03
$ws.Range("B41").Value = 'lengthKindPattern_03'
$ws.Range("C41").Value = 'Passed'
$ws.Range("D41").Value = 'DFDL-12-088R'
$ws.Range("E41").Value = 'High'
$ws.Range("F41").Value = 3
$ws.Range("G41").Value = 'text'
$ws.Range("H41").Value = 'DFDL-205'

# Row 42: litNil1
$ws.Range("B42").Value = 'litNil1'
$ws.Range("C42").Value = 'Passed'
$ws.Range("D42").Value = 'DFDL-13-234R'
$ws.Range("E42").Value = 'High'
$ws.Range("F42").Value = 3
$ws.Range("G42").Value = 'text'
$ws.Range("H42").Value = 'DFDL-199'

# Row 43: litNil2
$ws.Range("B43").Value = 'litNil2'
$ws.Range("C43").Value = 'Passed'
$ws.Range("D43").Value = 'DFDL-13-234R'
$ws.Range("E43").Value = 'High'
$ws.Range("F43").Value = 3
$ws.Range("G43").Value = 'text'
$ws.Range("H43").Value = 'DFDL-199'

# Row 44: multiple_delimiters2
$ws.Range("B44").Value = 'multiple_delimiters2'
$ws.Range("C44").Value = 'Passed'
$ws.Range("D44").Value = 'DFDL-14-008R'
$ws.Range("E44").Value = 'High'
$ws.Range("F44").Value = 3
$ws.Range("G44").Value = 'text'
$ws.Range("H44").Value = 'DFDL-109'

# Row 45: basic
$ws.Range("B45").Value = 'basic'
$ws.Range("C45").Value = 'Passed'
$ws.Range("D45").Value = 'DFDL-15-001R'
$ws.Range("E45").Value = 'High'
$ws.Range("F45").Value = 3
$ws.Range("G45").Value = 'text'
$ws.Range("H45").Value = 'DFDL-204'

# Row 46: choice2
$ws.Range("B46").Value = 'choice2'
$ws.Range("C46").Value = 'Passed'
$ws.Range("D46").Value = 'DFDL-15-001R'
$ws.Range("E46").Value = 'High'
$ws.Range("F46").Value = 3
$ws.Range("G46").Value = 'text'
$ws.Range("H46").Value = 'DFDL-204'

# Row 47: choice3
$ws.Range("B47").Value = 'choice3'
$ws.Range("C47").Value = 'Passed'
$ws.Range("D47").Value = 'DFDL-15-001R'
$ws.Range("E47").Value = 'High'
$ws.Range("F47").Value = 3
$ws.Range("G47").Value = 'text'
$ws.Range("H47").Value = 'DFDL-204'

# Row 48: choice4
$ws.Range("B48").Value = 'choice4'
$ws.Range("C48").Value = 'Passed'
$ws.Range("D48").Value = 'DFDL-15-001R'
$ws.Range("E48").Value = 'High'
$ws.Range("F48").Value = 3
$ws.Range("G48").Value = 'text'
$ws.Range("H48").Value = 'DFDL-204'

# Row 49: choice5
$ws.Range("B49").Value = 'choice5'
$ws.Range("C49").Value = 'Passed'
$ws.Range("D49").Value = 'DFDL-15-001R'
$ws.Range("E49").Value = 'High'
$ws.Range("F49").Value = 3
$ws.Range("G49").Value = 'text'
$ws.Range("H49").Value = 'DFDL-204'

# Row 50: choice6
$ws.Range("B50").Value = 'choice6'
$ws.Range("C50").Value = 'Passed'
$ws.Range("D50").Value = 'DFDL-15-001R'
$ws.Range("E50").Value = 'High'
$ws.Range("F50").Value = 3
$ws.Range("G50").Value = 'text'
$ws.Range("H50").Value = 'DFDL-204'

# Row 51: choiceFail1
$ws.Range("B51").Value = 'choiceFail1'
$ws.Range("C51").Value = 'Passed'
$ws.Range("D51").Value = 'DFDL-15-001R'
$ws.Range("E51").Value = 'High'
$ws.Range("F51").Value = 3
$ws.Range("G51").Value = 'text'
$ws.Range("H51").Value = 'DFDL-204'
$ws.Range("I51").Value = 'Yes'

# Row 52: choiceDelim1
$ws.Range("B52").Value = 'choiceDelim1'
$ws.Range("C52").Value = 'Passed'
$ws.Range("D52").Value = 'DFDL-15-001R'
$ws.Range("E52").Value = 'High'
$ws.Range("F52").Value = 3
$ws.Range("G52").Value = 'text'
$ws.Range("H52").Value = 'DFDL-204'

# Row 53: nestedChoice1
$ws.Range("B53").Value = 'nestedChoice1'
$ws.Range("C53").Value = 'Passed'
$ws.Range("D53").Value = 'DFDL-15-001R'
$ws.Range("E53").Value = 'High'
$ws.Range("F53").Value = 3
$ws.Range("G53").Value = 'text'
$ws.Range("H53").Value = 'DFDL-204'

# Row 54: property_scoping_01
$ws.Range("B54").Value = 'property_scoping_01'
$ws.Range("C54").Value = 'Passed'
$ws.Range("D54").Value = 'DFDL-8-009R'
$ws.Range("E54").Value = 'High'
$ws.Range("F54").Value = 3
$ws.Range("G54").Value = 'text'
$ws.Range("H54").Value = 'DFDL-224'

# Row 55: syntax_entities_6_02
$ws.Range("B55").Value = 'syntax_entities_6_02'
$ws.Range("C55").Value = 'Passed'
$ws.Range("D55").Value = 'DFDL-6-042R'
$ws.Range("E55").Value = 'High'
$ws.Range("F55").Value = 3
$ws.Range("G55").Value = 'byte'
$ws.Range("H55").Value = 'DFDL-219'

# Final selection matches the saved workbook state (last touched cell)
$ws.Range("G55").Select() | Out-Null
